# Update Activity List, Gantt Chart, WBS
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prepare the two new rows (18, 19) by copying the formatting of an
# existing fully-populated data row (row 16: styles 1,1,1,1,4,1,1) so the
# new cells get the right number formats (in particular the date style on
# column E) without inventing new style/numFmt entries.
$ws.Range("A16:G16").Copy() | Out-Null
$ws.Range("A18:G19").PasteSpecial(-4122) | Out-Null
$ws.Range("A18:G19").ClearContents()

# --- Write new shared-string values in the exact order they first appear
# so the regenerated shared-strings table lines up with the target.
$ws.Cells.Item(16, 3).Value = "Business Rules"          # C16
$ws.Cells.Item(17, 3).Value = "ERD"                      # C17
$ws.Cells.Item(18, 2).Value = "Development"              # B18
$ws.Cells.Item(18, 3).Value = "Admin and User Interfaces"# C18
$ws.Cells.Item(4, 6).Value  = "2;3"                       # F4
$ws.Cells.Item(6, 6).Value  = "5;6"                       # F6
$ws.Cells.Item(7, 6).Value  = "5;6;7"                     # F7
$ws.Cells.Item(8, 6).Value  = "5;6;7;8"                   # F8
$ws.Cells.Item(9, 6).Value  = "5;9"                       # F9
$ws.Cells.Item(13, 6).Value = "3;5;6;7;8"                 # F13
$ws.Cells.Item(15, 6).Value = "3;5;6;7;8;9;10;15"         # F15
$ws.Cells.Item(16, 6).Value = "5;6;7;8;9;10"               # F16
$ws.Cells.Item(18, 7).Value = "Develop Interfaces"        # G18
$ws.Cells.Item(15, 7).Value = "Updated system "           # G15
$ws.Cells.Item(14, 7).Value = "Yii Framework"             # G14

# --- Remaining value changes (rows 3-17), values only, order irrelevant
# since either numeric or already-known shared strings.
$ws.Cells.Item(3, 6).Value = 2                            # F3
$ws.Cells.Item(3, 7).Value = "Plan for project features"  # G3

$ws.Cells.Item(4, 7).Value = "Additional Information"     # G4

$ws.Cells.Item(5, 6).Value = 5                             # F5
$ws.Cells.Item(5, 7).Value = "Finalized Event Table"      # G5

$ws.Cells.Item(6, 7).Value = "Finalized Use Case Full Description" # G6

$ws.Cells.Item(7, 7).Value = "Finalized all UML Diagrams" # G7

$ws.Cells.Item(8, 7).Value = "Updated Wiki and OneNote "   # G8

$ws.Cells.Item(9, 1).Value = 9                              # A9
$ws.Cells.Item(9, 3).Value = "Project Requirements"        # C9
$ws.Cells.Item(9, 4).Value = 6                              # D9
$ws.Cells.Item(9, 5).Value = 42650                          # E9
$ws.Cells.Item(9, 7).Value = "Prepared all project requirements" # G9

$ws.Cells.Item(10, 1).Value = 10                            # A10
$ws.Cells.Item(10, 3).Value = "Activity List"               # C10
$ws.Cells.Item(10, 4).Value = 89                             # D10
$ws.Cells.Item(10, 5).Value = 42548                          # E10
$ws.Cells.Item(10, 6).Clear()                                 # F10 (removed)
$ws.Cells.Item(10, 7).Value = "Activity List"                # G10

$ws.Cells.Item(11, 1).Value = 12                            # A11
$ws.Cells.Item(11, 3).Value = "Develop WBS"                 # C11
$ws.Cells.Item(11, 4).Value = 4                              # D11
$ws.Cells.Item(11, 5).Value = 42663                          # E11
$ws.Cells.Item(11, 6).Value = 11                              # F11 (added)
$ws.Cells.Item(11, 7).Value = "WBS"                          # G11

$ws.Cells.Item(12, 1).Value = 13                            # A12
$ws.Cells.Item(12, 3).Value = "Develop Gantt Chart"          # C12
$ws.Cells.Item(12, 4).Value = 4                               # D12
$ws.Cells.Item(12, 5).Value = 42666                           # E12
$ws.Cells.Item(12, 6).Value = 11                               # F12 (added)
$ws.Cells.Item(12, 7).Value = "Gantt Chart"                   # G12

$ws.Cells.Item(13, 1).Value = 8                              # A13
$ws.Cells.Item(13, 2).Value = "Design"                       # B13
$ws.Cells.Item(13, 3).Value = "Develop GUI"                  # C13
$ws.Cells.Item(13, 4).Value = 14                              # D13
$ws.Cells.Item(13, 5).Value = 42608                           # E13
$ws.Cells.Item(13, 7).Value = "GUI of the proposed system"   # G13

$ws.Cells.Item(14, 1).Value = 14                             # A14
$ws.Cells.Item(14, 2).Value = "Design"                       # B14
$ws.Cells.Item(14, 3).Value = "Install Yii"                  # C14
$ws.Cells.Item(14, 5).Value = 42583                           # E14
$ws.Cells.Item(14, 6).Clear()                                 # F14 (removed)

$ws.Cells.Item(15, 1).Value = 15                             # A15
$ws.Cells.Item(15, 3).Value = "Design System"                # C15
$ws.Cells.Item(15, 4).Value = 60                               # D15
$ws.Cells.Item(15, 5).Value = 42628                            # E15

$ws.Cells.Item(16, 1).Value = 16                             # A16
$ws.Cells.Item(16, 4).Value = 8                                # D16
$ws.Cells.Item(16, 5).Value = 42638                            # E16
$ws.Cells.Item(16, 7).Value = "Business Rules"                # G16

$ws.Cells.Item(17, 1).Clear()                                 # A17 (removed entirely)
$ws.Cells.Item(17, 4).Value = 21                               # D17
$ws.Cells.Item(17, 6).Value = 5                                 # F17
$ws.Cells.Item(17, 7).Value = "ERD"                             # G17

# --- New row 18
$ws.Cells.Item(18, 1).Value = 11                             # A18
$ws.Cells.Item(18, 4).Value = 60                               # D18
$ws.Cells.Item(18, 5).Value = 42660                            # E18
$ws.Cells.Item(18, 6).Value = "5;6;7;8;9;10"                   # F18

# --- New row 19
$ws.Cells.Item(19, 6).Clear()                                 # F19 stays empty (no cell at all)
$ws.Cells.Item(19, 1).Value = 11                             # A19
$ws.Cells.Item(19, 2).Value = "Development"                  # B19
$ws.Cells.Item(19, 3).Value = "Bluemix"                       # C19
$ws.Cells.Item(19, 4).Value = 11                               # D19
$ws.Cells.Item(19, 5).Value = 42636                            # E19
$ws.Cells.Item(19, 7).Value = "Application for cloud hosting" # G19

# --- View state: scroll / selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("G15").Select()
